$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.065.20"
$ws.Range("D3").Value = "3.366.58"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "3.946.40"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.364.76"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "61.138.74"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "3.504.83"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000108"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.84%  "
$ws.Range("E27").Value = "  -3.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "170.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -4.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "29.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.96%  "
$ws.Range("D40").Value = "3.404.57"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0754"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.761"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("D46").Value = "2.488.36"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.00%  "
